# Update cryptocurrency price (column D) and volume change (column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.478.23"
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = "'1.864.86"
$ws.Range("E3").Value = '  +2.70%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = "'316.02"
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = "'0.4661"
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").Value = "'0.3731"
$ws.Range("E8").Value = '  +2.15%  '
$ws.Range("D9").Value = "'0.07385"
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("D10").Value = "'0.8880"
$ws.Range("E10").Value = '  +3.56%  '
$ws.Range("D11").Value = "'0.07974"
$ws.Range("E11").Value = '  +5.84%  '
$ws.Range("D12").Value = "'19.99"
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = "'1.867.28"
$ws.Range("E13").Value = '  +4.83%  '
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("D15").Value = "'6.593"
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = "'0.000008963"
$ws.Range("E18").Value = '  +3.99%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = "'14.94"
$ws.Range("E20").Value = '  +3.54%  '
$ws.Range("D21").Value = "'27.498.13"
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("D22").Value = "'5.175"
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = "'2.079.72"
$ws.Range("E24").Value = '  +6.54%  '
$ws.Range("D25").Value = "'153.29"
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").Value = "'1.882"
$ws.Range("E26").Value = '  +2.53%  '
$ws.Range("D27").Value = "'18.57"
$ws.Range("E27").Value = '  +2.60%  '
$ws.Range("D28").Value = "'2.088"
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").Value = "'5.165"
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("D30").Value = "'117.17"
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("D31").Value = "'0.08910"
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("D32").Value = "'3.012"
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("E33").Value = '  +5.52%  '
$ws.Range("D34").Value = "'1.161"
$ws.Range("E34").Value = '  +3.02%  '
$ws.Range("D35").Value = "'4.485"
$ws.Range("E35").Value = '  +1.82%  '
$ws.Range("D36").Value = "'2.648"
$ws.Range("E36").Value = '  +9.47%  '
$ws.Range("D37").Value = "'0.01971"
$ws.Range("E37").Value = '  +2.98%  '
$ws.Range("D38").Value = "'1.083"
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").Value = "'0.05291"
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D40").Value = "'2.987"
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").Value = "'7.185"
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("D42").Value = "'0.5218"
$ws.Range("D43").Value = "'0.1647"
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").Value = "'8.363"
$ws.Range("E44").Value = '  +2.58%  '
$ws.Range("D45").Value = "'0.4891"
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("D46").Value = "'10.30"
$ws.Range("E46").Value = '  +2.50%  '
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").Value = "'1.665"
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("D49").Value = "'103.48"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").Value = "'0.06260"
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").Value = "'65.77"
$ws.Range("E51").Value = '  +3.06%  '
